$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    701.3070661196912,
    701.3070661196912,
    743.3819643578375,
    828.8832227764427,
    916.252090899261,
    1005.571982782548,
    1096.932028729645,
    1190.427609826836,
    1286.160956459582,
    1384.241820224659,
    828.8832227764427,
    828.8832227764427,
    828.8832227764427,
    828.8832227764427,
    916.252090899261,
    1005.571982782548,
    1096.932028729645,
    1190.427609826836,
    1286.160956459582,
    1384.241820224659,
    828.8832227764427,
    828.8832227764427,
    828.8832227764427,
    828.8832227764427,
    916.252090899261,
    1005.571982782548,
    1096.932028729645,
    1190.427609826836,
    1286.160956459582,
    1384.241820224659
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}
